$d = $word.ActiveDocument

# Update the date heading (unique text, safe to use Find/Replace).
$d.Content.Find.Execute("2024-08-11 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-08-12 Monday", 2)

# Update the division-problem table cells. Several cells share identical
# text (e.g. "34÷4=8, 2" appears twice), so address cells by their
# (row, column) coordinates instead of a global text search/replace.
$t = $word.ActiveDocument.Tables.Item(1)

$updates = @(
    @(1, 1, "27÷9=3, 0"),
    @(1, 3, "99÷2=49, 1"),
    @(1, 4, "50÷2=25, 0"),
    @(1, 5, "62÷7=8, 6"),
    @(5, 1, "29÷8=3, 5"),
    @(5, 2, "59÷4=14, 3"),
    @(5, 3, "22÷8=2, 6"),
    @(5, 4, "43÷7=6, 1"),
    @(5, 5, "54÷7=7, 5"),
    @(9, 1, "78÷6=13, 0"),
    @(9, 2, "76÷6=12, 4"),
    @(9, 3, "74÷2=37, 0"),
    @(9, 4, "91÷6=15, 1"),
    @(9, 5, "79÷7=11, 2"),
    @(13, 1, "21÷9=2, 3"),
    @(13, 2, "15÷8=1, 7"),
    @(13, 3, "58÷7=8, 2"),
    @(13, 4, "64÷4=16, 0"),
    @(13, 5, "48÷6=8, 0"),
    @(17, 1, "77÷6=12, 5"),
    @(17, 2, "82÷4=20, 2"),
    @(17, 3, "92÷9=10, 2"),
    @(17, 4, "94÷6=15, 4"),
    @(17, 5, "58÷9=6, 4")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $newText = $u[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
